# Updated symbol list on Sun Feb  5 10:29:10 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) for the coin rows
# that changed. Values are written with a leading apostrophe so Excel
# keeps them as literal text (matching the workbook's existing
# inlineStr/text cell format) instead of auto-converting to numbers or
# percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'334.91"
$ws.Range("E2").Value = "'1.69%"
$ws.Range("D3").Value = "'43.88"
$ws.Range("E3").Value = "'6.17%"
$ws.Range("D4").Value = "'5.738"
$ws.Range("E4").Value = "'2.22%"
$ws.Range("D5").Value = "'0.08339"
$ws.Range("E5").Value = "'1.77%"
$ws.Range("D6").Value = "'8.834"
$ws.Range("E6").Value = "'0.92%"
$ws.Range("D7").Value = "'1.965"
$ws.Range("E7").Value = "'-1.98%"
$ws.Range("E8").Value = "'-3.60%"
$ws.Range("D9").Value = "'0.9438"
$ws.Range("E9").Value = "'2.27%"
$ws.Range("D10").Value = "'0.1247"
$ws.Range("E10").Value = "'-2.39%"
$ws.Range("D11").Value = "'0.1983"
$ws.Range("E11").Value = "'1.56%"
$ws.Range("D12").Value = "'0.1073"
$ws.Range("E12").Value = "'15.70%"
$ws.Range("D13").Value = "'0.04537"
$ws.Range("E13").Value = "'18.08%"
$ws.Range("E14").Value = "'0.81%"
$ws.Range("D15").Value = "'0.001296"
$ws.Range("E15").Value = "'-0.84%"
$ws.Range("D16").Value = "'0.005950"
$ws.Range("E16").Value = "'-4.47%"
$ws.Range("E17").Value = "'1.54%"
$ws.Range("D18").Value = "'4.518"
$ws.Range("E18").Value = "'0.42%"
$ws.Range("E19").Value = "'0.74%"
$ws.Range("D20").Value = "'8.711"
$ws.Range("E20").Value = "'5.93%"
$ws.Range("D21").Value = "'0.1352"
$ws.Range("E21").Value = "'-0.93%"
$ws.Range("D22").Value = "'0.2690"
$ws.Range("E22").Value = "'1.16%"
$ws.Range("D23").Value = "'0.04411"
$ws.Range("E23").Value = "'0.07%"
$ws.Range("D24").Value = "'0.001255"
$ws.Range("E24").Value = "'-0.15%"
$ws.Range("D25").Value = "'0.004363"
$ws.Range("E25").Value = "'1.08%"
$ws.Range("E26").Value = "'5.08%"
$ws.Range("D39").Value = "'0.02814"
$ws.Range("E39").Value = "'2.32%"
$ws.Range("D40").Value = "'0.06030"
$ws.Range("E40").Value = "'10.49%"
$ws.Range("D41").Value = "'0.007932"
$ws.Range("E41").Value = "'1.67%"
$ws.Range("D42").Value = "'0.1428"
$ws.Range("E42").Value = "'0.46%"
$ws.Range("D43").Value = "'0.008966"
$ws.Range("E43").Value = "'0.33%"
$ws.Range("D44").Value = "'0.002172"
$ws.Range("E44").Value = "'0.03%"
$ws.Range("D45").Value = "'0.01015"
$ws.Range("E45").Value = "'-11.29%"
$ws.Range("E46").Value = "'3.48%"
$ws.Range("E47").Value = "'0.08%"
$ws.Range("D48").Value = "'0.003188"
$ws.Range("E48").Value = "'-0.09%"
$ws.Range("E49").Value = "'-0.32%"
$ws.Range("E50").Value = "'0.08%"
$ws.Range("E51").Value = "'0.08%"
